$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 7318.478
$ws.Range("I98").Value = 5326.25
$ws.Range("J98").Value = 20600
$ws.Range("K98").Value = 5326.25
$ws.Range("L98").Value = 20600
$ws.Range("M98").Value = -3828.25
$ws.Range("N98").Value = -23596

$ws.Range("H122").Value = 7318.478
$ws.Range("I122").Value = 5326.25
$ws.Range("J122").Value = 20600
$ws.Range("K122").Value = 15978.75
$ws.Range("L122").Value = 61800
$ws.Range("M122").Value = -13528.75
$ws.Range("N122").Value = -66700

$ws.Range("H137").Value = 1163.9215
$ws.Range("I137").Value = 1126.3636
$ws.Range("J137").Value = 1232.7778
$ws.Range("K137").Value = 3379.0908
$ws.Range("L137").Value = 3698.3334
$ws.Range("M137").Value = -829.0907999999999
$ws.Range("N137").Value = -8798.3334

$ws.Range("H140").Value = 84532.5
$ws.Range("I140").Value = 37500
$ws.Range("J140").Value = 91251.42999999999
$ws.Range("K140").Value = 37500
$ws.Range("L140").Value = 91251.42999999999
$ws.Range("M140").Value = -32320
$ws.Range("N140").Value = -101611.43

$ws.Range("H141").Value = 7778.7144
$ws.Range("I141").Value = 4492.885
$ws.Range("J141").Value = 50494.5
$ws.Range("K141").Value = 13478.655
$ws.Range("L141").Value = 151483.5
$ws.Range("M141").Value = -8298.655000000001
$ws.Range("N141").Value = -161843.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1485.3334
$ws.Range("I45").Value = 1432.4
$ws.Range("K45").Value = 1432.4
$ws.Range("M45").Value = -1055.4

$ws.Range("H53").Value = 35000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 35000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 35000
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -36364

$ws.Range("H61").Value = 1648.8
$ws.Range("I61").Value = 1382.1177
$ws.Range("J61").Value = 3160
$ws.Range("K61").Value = 1382.1177
$ws.Range("L61").Value = 3160
$ws.Range("M61").Value = -1170.1177
$ws.Range("N61").Value = -3584

$ws.Range("H74").Value = 1393.1666
$ws.Range("I74").Value = 1307.091
$ws.Range("K74").Value = 1307.091
$ws.Range("M74").Value = -433.0909999999999

$ws.Range("H77").Value = 1393.1666
$ws.Range("I77").Value = 1307.091
$ws.Range("K77").Value = 6535.455
$ws.Range("M77").Value = -2167.455

$ws.Range("H122").Value = 1946.9678
$ws.Range("I122").Value = 1909.5927
$ws.Range("J122").Value = 2199.25
$ws.Range("K122").Value = 5728.7781
$ws.Range("L122").Value = 6597.75
$ws.Range("M122").Value = -3278.7781
$ws.Range("N122").Value = -11497.75

$ws.Range("H132").Value = 835050.3
$ws.Range("I132").Value = 1429263.9
$ws.Range("J132").Value = 3151.2
$ws.Range("K132").Value = 4287791.699999999
$ws.Range("L132").Value = 9453.599999999999
$ws.Range("M132").Value = -4285261.699999999
$ws.Range("N132").Value = -14513.6

$ws.Range("H136").Value = 1648.8
$ws.Range("I136").Value = 1382.1177
$ws.Range("J136").Value = 3160
$ws.Range("K136").Value = 4146.3531
$ws.Range("L136").Value = 9480
$ws.Range("M136").Value = -1596.3531
$ws.Range("N136").Value = -14580

$ws.Range("H141").Value = 46273.6
$ws.Range("J141").Value = 46273.6
$ws.Range("L141").Value = 46273.6
$ws.Range("N141").Value = -56633.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 978612.75
$ws.Range("I134").Value = 1748054.4
$ws.Range("J134").Value = 3986.6667
$ws.Range("K134").Value = 5244163.199999999
$ws.Range("L134").Value = 11960.0001
$ws.Range("M134").Value = -5241628.199999999
$ws.Range("N134").Value = -17030.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1248.5
$ws.Range("I16").Value = 1350
$ws.Range("K16").Value = 1350
$ws.Range("M16").Value = -1063

$ws.Range("H28").Value = 197428.33
$ws.Range("J28").Value = 197428.33
$ws.Range("L28").Value = 197428.33
$ws.Range("N28").Value = -197918.33

$ws.Range("H31").Value = 12822792
$ws.Range("I31").Value = 19609358
$ws.Range("J31").Value = 3724.8518
$ws.Range("K31").Value = 19609358
$ws.Range("L31").Value = 3724.8518
$ws.Range("M31").Value = -19609063
$ws.Range("N31").Value = -4314.8518

$ws.Range("H34").Value = 12822792
$ws.Range("I34").Value = 19609358
$ws.Range("J34").Value = 3724.8518
$ws.Range("K34").Value = 19609358
$ws.Range("L34").Value = 3724.8518
$ws.Range("M34").Value = -19609156
$ws.Range("N34").Value = -4128.8518

$ws.Range("H58").Value = 1412.4242
$ws.Range("I58").Value = 1393.24
$ws.Range("J58").Value = 1472.375
$ws.Range("K58").Value = 1393.24
$ws.Range("L58").Value = 1472.375
$ws.Range("M58").Value = -1190.24
$ws.Range("N58").Value = -1878.375

$ws.Range("H113").Value = 1248.5
$ws.Range("I113").Value = 1350
$ws.Range("K113").Value = 1350
$ws.Range("M113").Value = 820

$ws.Range("H122").Value = 1003.9286
$ws.Range("I122").Value = 988.5217
$ws.Range("K122").Value = 2965.5651
$ws.Range("M122").Value = -515.5650999999998

$ws.Range("H132").Value = 2258.1462
$ws.Range("I132").Value = 1987.7428
$ws.Range("J132").Value = 3835.5
$ws.Range("K132").Value = 5963.2284
$ws.Range("L132").Value = 11506.5
$ws.Range("M132").Value = -3433.2284
$ws.Range("N132").Value = -16566.5

$ws.Range("H134").Value = 359977.84
$ws.Range("I134").Value = 488959.16
$ws.Range("J134").Value = 1696.4445
$ws.Range("K134").Value = 1466877.48
$ws.Range("L134").Value = 5089.333500000001
$ws.Range("M134").Value = -1464342.48
$ws.Range("N134").Value = -10159.3335

$ws.Range("H135").Value = 80590.414
$ws.Range("J135").Value = 80590.414
$ws.Range("L135").Value = 80590.414
$ws.Range("N135").Value = -90730.414

$ws.Range("H136").Value = 1412.4242
$ws.Range("I136").Value = 1393.24
$ws.Range("J136").Value = 1472.375
$ws.Range("K136").Value = 4179.72
$ws.Range("L136").Value = 4417.125
$ws.Range("M136").Value = -1629.72
$ws.Range("N136").Value = -9517.125

$ws.Range("H141").Value = 23123.25
$ws.Range("J141").Value = 23123.25
$ws.Range("L141").Value = 23123.25
$ws.Range("N141").Value = -33483.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 966238.6
$ws.Range("J12").Value = 1136726
$ws.Range("L12").Value = 3410178
$ws.Range("N12").Value = -3410524

$ws.Range("H80").Value = 5283.25
$ws.Range("J80").Value = 3139.9
$ws.Range("L80").Value = 9419.700000000001
$ws.Range("N80").Value = -11291.7

$ws.Range("H83").Value = 5283.25
$ws.Range("J83").Value = 3139.9
$ws.Range("L83").Value = 28259.1
$ws.Range("N83").Value = -37619.10000000001

$ws.Range("H131").Value = 859.76
$ws.Range("J131").Value = 883.8936
$ws.Range("L131").Value = 2651.6808
$ws.Range("N131").Value = -12731.6808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 19980
$ws.Range("J24").Value = 19980
$ws.Range("L24").Value = 19980
$ws.Range("N24").Value = -20326

$ws.Range("H101").Value = 60492.09
$ws.Range("J101").Value = 60492.09
$ws.Range("L101").Value = 60492.09
$ws.Range("N101").Value = -66982.09

$ws.Range("H102").Value = 2104.9062
$ws.Range("I102").Value = 2079.087
$ws.Range("J102").Value = 2170.889
$ws.Range("K102").Value = 2079.087
$ws.Range("L102").Value = 2170.889
$ws.Range("M102").Value = -457.087
$ws.Range("N102").Value = -5414.889

$ws.Range("H109").Value = 17811.545
$ws.Range("J109").Value = 17811.545
$ws.Range("L109").Value = 17811.545
$ws.Range("N109").Value = -19891.545

$ws.Range("H113").Value = 1086.6666
$ws.Range("I113").Value = 966.6923
$ws.Range("J113").Value = 1398.6
$ws.Range("K113").Value = 966.6923
$ws.Range("L113").Value = 1398.6
$ws.Range("M113").Value = 1203.3077
$ws.Range("N113").Value = -5738.6

$ws.Range("H122").Value = 2630.12
$ws.Range("I122").Value = 2673.366
$ws.Range("J122").Value = 2433.111
$ws.Range("K122").Value = 8020.098
$ws.Range("L122").Value = 7299.333
$ws.Range("M122").Value = -5570.098
$ws.Range("N122").Value = -12199.333

$ws.Range("H123").Value = 17963.916
$ws.Range("J123").Value = 17963.916
$ws.Range("L123").Value = 17963.916
$ws.Range("N123").Value = -22863.916

$ws.Range("H132").Value = 2153.6099
$ws.Range("I132").Value = 2022.5
$ws.Range("J132").Value = 2338.7058
$ws.Range("K132").Value = 6067.5
$ws.Range("L132").Value = 7016.117400000001
$ws.Range("M132").Value = -3537.5
$ws.Range("N132").Value = -12076.1174

$ws.Range("H133").Value = 56856.25
$ws.Range("J133").Value = 56856.25
$ws.Range("L133").Value = 56856.25
$ws.Range("N133").Value = -66976.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H86").Value = 31195
$ws.Range("J86").Value = 31195
$ws.Range("L86").Value = 31195
$ws.Range("N86").Value = -33567

$ws.Range("H89").Value = 31195
$ws.Range("J89").Value = 31195
$ws.Range("L89").Value = 93585
$ws.Range("N89").Value = -105441

$ws.Range("H136").Value = 4198.6484
$ws.Range("I136").Value = 4358.273
$ws.Range("J136").Value = 2881.75
$ws.Range("K136").Value = 13074.819
$ws.Range("L136").Value = 8645.25
$ws.Range("M136").Value = -10524.819
$ws.Range("N136").Value = -13745.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 48966
$ws.Range("J68").Value = 48966
$ws.Range("L68").Value = 48966
$ws.Range("N68").Value = -50588

$ws.Range("H71").Value = 48966
$ws.Range("J71").Value = 48966
$ws.Range("L71").Value = 146898
$ws.Range("N71").Value = -155010

$ws.Range("H122").Value = 25003698
$ws.Range("I122").Value = 35716428
$ws.Range("K122").Value = 107149284
$ws.Range("M122").Value = -107146834

$ws.Range("H132").Value = 2704.9614
$ws.Range("I132").Value = 2029.6875
$ws.Range("J132").Value = 3785.4
$ws.Range("K132").Value = 6089.0625
$ws.Range("L132").Value = 11356.2
$ws.Range("M132").Value = -3559.0625
$ws.Range("N132").Value = -16416.2

$ws.Range("H136").Value = 1572.3077
$ws.Range("I136").Value = 1640.4054
$ws.Range("K136").Value = 4921.216200000001
$ws.Range("M136").Value = -2371.216200000001
